$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 2.75
$ws.Range("H2").Value = 3.6
$ws.Range("I2").Value = 2.4
$ws.Range("J2").Value = 3.25
$ws.Range("K2").Value = 2.25
$ws.Range("L2").Value = 3
$ws.Range("M2").Value = 1.04
$ws.Range("N2").Value = 13
$ws.Range("O2").Value = 1.25
$ws.Range("P2").Value = 4
$ws.Range("Q2").Value = 1.73
$ws.Range("R2").Value = 2.1
$ws.Range("S2").Value = 1.33
$ws.Range("T2").Value = 3.25
$ws.Range("U2").Value = 1.62
$ws.Range("V2").Value = 2.2
$ws.Range("W2").Value = 10
$ws.Range("X2").Value = 15
$ws.Range("Y2").Value = 10
$ws.Range("AB2").Value = 26
$ws.Range("AC2").Value = 13
$ws.Range("AD2").Value = 7
$ws.Range("AE2").Value = 13
$ws.Range("AG2").Value = 151
$ws.Range("AH2").Value = 9.5
$ws.Range("AK2").Value = 23
$ws.Range("AM2").Value = 23
$ws.Range("AR2").Value = 67
$ws.Range("AT2").Value = 3.25
$ws.Range("AU2").Value = 7.5
$ws.Range("AW2").Value = 501
$ws.Range("AY2").Value = 13
$ws.Range("AZ2").Value = 21
$ws.Range("BC2").Value = 126
$ws.Range("BD2").Value = 126

# Row 3
$ws.Range("G3").Value = 2.7
$ws.Range("H3").Value = 3.2
$ws.Range("I3").Value = 2.45
$ws.Range("J3").Value = 3.5
$ws.Range("L3").Value = 3.25
$ws.Range("M3").Value = 1.06
$ws.Range("N3").Value = 10
$ws.Range("U3").Value = 1.83
$ws.Range("V3").Value = 1.83
$ws.Range("Y3").Value = 11
$ws.Range("Z3").Value = 29
$ws.Range("AA3").Value = 23
$ws.Range("AH3").Value = 8
$ws.Range("AI3").Value = 12
$ws.Range("AK3").Value = 23
$ws.Range("AN3").Value = 4.75
$ws.Range("AO3").Value = 17
$ws.Range("BB3").Value = 67

# Row 4
$ws.Range("G4").Value = 2.15
$ws.Range("H4").Value = 3.25
$ws.Range("J4").Value = 2.88
$ws.Range("K4").Value = 2.2
$ws.Range("M4").Value = 1.05
$ws.Range("N4").Value = 11
$ws.Range("O4").Value = 1.25
$ws.Range("P4").Value = 3.75
$ws.Range("Q4").Value = 1.9
$ws.Range("R4").Value = 1.95
$ws.Range("W4").Value = 8.5
$ws.Range("Y4").Value = 9
$ws.Range("AA4").Value = 17
$ws.Range("AB4").Value = 26
$ws.Range("AC4").Value = 11
$ws.Range("AD4").Value = 6.5
$ws.Range("AF4").Value = 41
$ws.Range("AI4").Value = 17
$ws.Range("AO4").Value = 12
$ws.Range("AP4").Value = 21
$ws.Range("AR4").Value = 51
$ws.Range("AU4").Value = 7.5
$ws.Range("BC4").Value = 151
